$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("responses")

# Update header row (row 1) labels for columns F..K
$ws.Range("F1").Value = "Final Marks"
$ws.Range("G1").Value = "Mathamatics"
$ws.Range("H1").Value = "Reasoning"
$ws.Range("I1").Value = "English"
$ws.Range("J1").Value = "GK"
$ws.Range("K1").Value = "Computer"

# Update data row (row 2) final marks value
$ws.Range("F2").Value = 341

# Remove the now-unused columns L through AE (headers + data)
$ws.Range("L1:AE2").Clear()
